# Insert a new data row above row 64 (pushing existing rows 64..160 down to
# 65..161) and populate it with the "Granada" price-record that the commit
# adds. This mirrors the XML diff, where row 64's old content reappears
# (unchanged) as the new row 65, and so on down to the old row 160 becoming
# the new row 161.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 64..160 down by one row.
$ws.Rows.Item(64).Insert()

# Populate the newly inserted row 64 with the new record's values.
$ws.Cells.Item(64, 1).Value = 10
$ws.Cells.Item(64, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(64, 3).Value = "La Araucanía"
$ws.Cells.Item(64, 4).Value = 44797
$ws.Cells.Item(64, 5).Value = 9
$ws.Cells.Item(64, 6).Value = "Fruta"
$ws.Cells.Item(64, 7).Value = 100104
$ws.Cells.Item(64, 8).Value = "Frutos de pepita"
$ws.Cells.Item(64, 9).Value = 100104001
$ws.Cells.Item(64, 10).Value = "Granada"
$ws.Cells.Item(64, 11).Value = "Wonderfull"
$ws.Cells.Item(64, 12).Value = "Primera"
$ws.Cells.Item(64, 13).Value = 35
$ws.Cells.Item(64, 14).Value = 14000
$ws.Cells.Item(64, 15).Value = 14000
$ws.Cells.Item(64, 16).Value = 14000
$ws.Cells.Item(64, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(64, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(64, 19).Value = 1400
$ws.Cells.Item(64, 20).Value = 10

# Match the date formatting used by the rest of column D.
$ws.Cells.Item(64, 4).NumberFormat = $ws.Cells.Item(65, 4).NumberFormat
